# Error Calculations and Plots
# Applies the missing-data re-roll: two rows ("RM 232" and "SC 92") are
# dropped entirely (remaining rows shift up), several cells flip between a
# present numeric value and a missing (blank) placeholder.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the "RM 232" row (original row 26) and the "SC 92" row
# (original row 28, which becomes row 27 once row 26 is gone). Everything
# below shifts up, turning the original A1:F35 range into A1:F33. ---
$ws.Rows.Item(26).Delete()
$ws.Rows.Item(27).Delete()

# --- Individual cell edits within the rows that kept their position (2-25) ---
$ws.Range("D3").Value = -14.2
$ws.Range("F4").ClearContents()
$ws.Range("D5").ClearContents()
$ws.Range("F9").Value = 17.26
$ws.Range("F10").Value = 16.43
$ws.Range("F11").Value = 17.65
$ws.Range("F12").Value = 17.45
$ws.Range("F15").ClearContents()
$ws.Range("F17").ClearContents()
$ws.Range("F18").ClearContents()
$ws.Range("F20").ClearContents()
$ws.Range("D21").Value = -14.3
$ws.Range("D23").ClearContents()

# --- Fill back in the values that were previously blanked in the rows that
# shifted up into 26-33 (SC 132's F and SC 193's D/F come back). ---
$ws.Range("F31").Value = 17.18
$ws.Range("D32").Value = -14.7
$ws.Range("F32").Value = 17.39
